# Add "Save" column (H) to the s_vals sheet, matching the header style of the
# existing columns and populating data rows with the save values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the neighboring header cell (G1) onto the new header
# cell H1 so it keeps the same bold/centered/bordered style used by the rest
# of row 1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Populate the new "Save" data values for each row.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
